# [pr4-2762] IPR: Create SPMetal configuration - synchronized with the
# current information model of the website.
#
# The "PCNCode" sheet's XML-mapped table (Tabela8) drops its "ProductName"
# column (originally column C): Title | ProductCodeNumber | ProductName |
# Compensation Good  ->  Title | ProductCodeNumber | Compensation Good.
# This also shrinks the shared-string table by the one now-unused entry,
# which is why every other sheet's string indices shift down by one.
# Finally, the PCNCode sheet becomes the active sheet/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PCNCode")
$lo = $ws.ListObjects.Item("Tabela8")

# Mirror the interactive selection (select column C) before removing it.
$ws.Columns.Item(3).Select()
$ws.Columns.Item(3).Delete()

# Shrink the table to the new 3-column extent and restore the header text
# for what is now column C ("Compensation Good", formerly column D).
$lo.Resize($ws.Range("A1:C5"))
$lo.ListColumns.Item(3).Range.Cells.Item(1, 1).Value = "Compensation Good"

# PCNCode ends up the active/selected sheet (activeTab moves from
# CutfillerCoefficient to PCNCode).
$ws.Activate()
